# Exporting excel - switch the D column from a duplicate birthday column
# into a BigDecimal "deposit" (存款) column, drop the old "年龄" (age)
# column E entirely, and rename the headers accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "年龄" column (E) completely so the sheet shrinks back to A:D.
$ws.Columns.Item(5).Delete()

# Header row: C1 loses the "转换测试" suffix, D1 becomes "存款".
$ws.Range("C1").Value = "生日"
$ws.Range("D1").Value = "存款"

# Replace column D's birthday-string duplicate with the numeric deposit values.
$ws.Range("D2").Value = 1.0
$ws.Range("D3").Value = 10.0
$ws.Range("D4").Value = 11.11
$ws.Range("D5").Value = 10.24
$ws.Range("D6").Value = 0.0
$ws.Range("D7").Value = 10.0
